$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.655.40"
$ws.Range("E2").Value = "  +3.53%  "

$ws.Range("D3").Value = "1.610.16"
$ws.Range("E3").Value = "  +2.88%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'212.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").Value = "'0.520"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.71%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").Value = "'27.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.67%  "

$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("E10").Value = "  +2.65%  "

$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("D12").Value = "'0.0910"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("D13").Value = "1.840.60"
$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("D14").Value = "1.656.27"
$ws.Range("E14").Value = "  +5.81%  "

$ws.Range("D15").Value = "29.667.45"
$ws.Range("E15").Value = "  +3.42%  "

$ws.Range("E16").Value = "  +4.29%  "

$ws.Range("E17").Value = "  +2.71%  "

$ws.Range("D18").Value = "'63.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.36%  "

$ws.Range("D19").Value = "'240.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.68%  "

$ws.Range("E20").Value = "  +3.89%  "

$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "'4.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.10%  "

$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("E26").Value = "  +2.28%  "

$ws.Range("E27").Value = "  +3.96%  "

$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("E29").Value = "  +2.92%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  +3.81%  "

$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").Value = "'3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("D34").Value = "'3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.89%  "

$ws.Range("D35").Value = "1.433.45"
$ws.Range("E35").Value = "  +1.95%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'1.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.95%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("E38").Value = "  +5.61%  "

$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").Value = "'0.0165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.25%  "

$ws.Range("D41").Value = "'0.540"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.46%  "

$ws.Range("E42").Value = "  +2.16%  "

$ws.Range("D43").Value = "'54.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +27.22%  "

$ws.Range("D44").Value = "'0.0491"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.03%  "

$ws.Range("D45").Value = "'0.802"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.63%  "

$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").Value = "'66.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.42%  "

$ws.Range("E48").Value = "  +1.60%  "

$ws.Range("D49").Value = "1.750.79"

$ws.Range("D50").Value = "'0.915"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.72%  "

$ws.Range("D51").Value = "'86.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.47%  "
